$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / rId1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 0
$ws1.Range("F3").Value = 108
$ws1.Range("F4").Value = 0
$ws1.Range("F5").Value = 0
$ws1.Range("F8").Value = 0
$ws1.Range("F10").Value = 463

# Sheet "全部类型" (sheet4 / rId4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 376
$ws4.Range("F3").Value = 108
$ws4.Range("F4").Value = 0
$ws4.Range("F5").Value = 14
$ws4.Range("F7").Value = 405
$ws4.Range("F8").Value = 140
$ws4.Range("F9").Value = 0
$ws4.Range("F10").Value = 0
